# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" quarterly sheet (cloned from the "2021-Q4"
#    sheet so it keeps the same column layout/styling) right before the
#    "总计" (total) sheet, and fill in the single fund holding row.
# 2. Prepend a new "2022-Q1" summary row to the "总计" sheet's table and
#    renumber the existing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: clone "2021-Q4" -> "2022-Q1", placed right after "2021-Q4"
# (i.e. right before "总计").
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)

$newQtr = $wb.ActiveSheet
$newQtr.Name = "2022-Q1"

# Overwrite the single data row with the new fund's info. Columns D:G
# are text-formatted numbers in the source data (keep leading/trailing
# zeros), column H is a real number.
$newQtr.Range("B2:G2").NumberFormat = "@"
$newQtr.Cells.Item(2, 2).Value = "001563"
$newQtr.Cells.Item(2, 3).Value = "华富健康文娱灵活配置混合"
$newQtr.Cells.Item(2, 4).Value = "0.10"
$newQtr.Cells.Item(2, 5).Value = "90.86"
$newQtr.Cells.Item(2, 6).Value = "2.82"
$newQtr.Cells.Item(2, 7).Value = "0.0028"
$newQtr.Cells.Item(2, 8).Value = 10
$newQtr.Range("B2:G2").Style = "Normal"

# ---------------------------------------------------------------------
# Step 2: update the "总计" sheet - insert a new top data row for
# 2022-Q1 and renumber the index column (A) for the rows that shift
# down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows("2:2").Insert()
$total.Range("B2:D2").Style = "Normal"

$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 1
$total.Cells.Item(2, 4).Value = 0

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4
